# TestRunner and excel sheet update
# Populates the "AddNewUserDetails" sheet (sheet2) with header row + sample
# user data row, including a text-stored phone number (quote-prefixed) and
# a hyperlink on the LinkedIn URL cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1)
$ws.Range("A1").Value = "firstname"
$ws.Range("B1").Value = "lastname"
$ws.Range("C1").Value = "location"
$ws.Range("D1").Value = "phone"
$ws.Range("E1").Value = "linkdeinurl"
$ws.Range("F1").Value = "usercomments"
$ws.Range("G1").Value = "postgraduate"
$ws.Range("H1").Value = "undergraduate"

# Data row (row 2)
$ws.Range("A2").Value = "Robert"
$ws.Range("B2").Value = "Davis"
$ws.Range("C2").Value = "Florida"
# Leading apostrophe forces text-with-quote-prefix storage (like typing
# '98675434 into a General formatted cell in Excel).
$ws.Range("D2").Value = "'98675434"
$ws.Range("E2").Value = "www.linkedin/ronj.com"
$ws.Range("F2").Value = "Testing"
$ws.Range("G2").Value = "Yes"
$ws.Range("H2").Value = "Yes"

# Hyperlink on the LinkedIn URL cell
$ws.Hyperlinks.Add($ws.Range("E2"), "http://www.linkedin/ronj.com")

# Select the used range, as in the saved workbook
[void]$ws.Range("A1:XFD2").Select()
